$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage across the affected range so Excel does not
# auto-coerce numeric-looking price/percentage strings into Number types.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '19.781.21'
$ws.Range("E2").Value = '  -8.74%  '
$ws.Range("D3").Value = '1.387.77'
$ws.Range("E3").Value = '  -9.53%  '
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.58%  '
$ws.Range("D5").Value = '1.005'
$ws.Range("E5").Value = '  +0.47%  '
$ws.Range("D6").Value = '269.03'
$ws.Range("E6").Value = '  -6.67%  '
$ws.Range("D7").Value = '0.3633'
$ws.Range("E7").Value = '  -7.61%  '
$ws.Range("D8").Value = '0.3018'
$ws.Range("E8").Value = '  -4.23%  '
$ws.Range("D9").Value = '38.59'
$ws.Range("E9").Value = '  -9.02%  '
$ws.Range("D10").Value = '0.06338'
$ws.Range("E10").Value = '  -11.38%  '
$ws.Range("D11").Value = '0.9585'
$ws.Range("E11").Value = '  -8.25%  '
$ws.Range("E12").Value = '  +0.73%  '
$ws.Range("D13").Value = '5.248'
$ws.Range("E13").Value = '  -7.00%  '
$ws.Range("D14").Value = '6.008'
$ws.Range("E14").Value = '  -8.80%  '
$ws.Range("D15").Value = '1.391.09'
$ws.Range("E15").Value = '  -9.43%  '
$ws.Range("D16").Value = '16.30'
$ws.Range("E16").Value = '  -12.03%  '
$ws.Range("D17").Value = '0.000009809'
$ws.Range("E17").Value = '  -9.96%  '
$ws.Range("D18").Value = '0.05631'
$ws.Range("E18").Value = '  -14.60%  '
$ws.Range("E19").Value = '  +0.50%  '
$ws.Range("D20").Value = '69.97'
$ws.Range("E20").Value = '  -15.96%  '
$ws.Range("D21").Value = '5.481'
$ws.Range("E21").Value = '  -10.31%  '
$ws.Range("D22").Value = '14.09'
$ws.Range("E22").Value = '  -8.55%  '
$ws.Range("D23").Value = '10.48'
$ws.Range("E23").Value = '  -2.98%  '
$ws.Range("D24").Value = '2.237'
$ws.Range("E24").Value = '  -5.20%  '
$ws.Range("D25").Value = '19.761.91'
$ws.Range("E25").Value = '  -8.84%  '
$ws.Range("D26").Value = '2.111'
$ws.Range("E26").Value = '  -9.91%  '
$ws.Range("D27").Value = '135.46'
$ws.Range("E27").Value = '  -8.47%  '
$ws.Range("D28").Value = '16.48'
$ws.Range("E28").Value = '  -10.01%  '
$ws.Range("D29").Value = '1.547.63'
$ws.Range("E29").Value = '  -9.45%  '
$ws.Range("D30").Value = '107.20'
$ws.Range("E30").Value = '  -8.30%  '
$ws.Range("D31").Value = '3.822'
$ws.Range("E31").Value = '  -20.96%  '
$ws.Range("D32").Value = '5.222'
$ws.Range("E32").Value = '  -11.19%  '
$ws.Range("D33").Value = '0.7897'
$ws.Range("E33").Value = '  -16.18%  '
$ws.Range("D34").Value = '0.07572'
$ws.Range("E34").Value = '  -6.87%  '
$ws.Range("D35").Value = '8.185'
$ws.Range("E35").Value = '  -3.80%  '
$ws.Range("D36").Value = '1.005'
$ws.Range("E36").Value = '  +0.47%  '
$ws.Range("D37").Value = '4.669'
$ws.Range("E37").Value = '  -8.76%  '
$ws.Range("D38").Value = '0.05529'
$ws.Range("E38").Value = '  -7.78%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '0.01995'
$ws.Range("E39").Value = '  -9.34%  '
$ws.Range("B40").Value = 'Algorand'
$ws.Range("C40").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D40").Value = '0.1870'
$ws.Range("E40").Value = '  -7.09%  '
$ws.Range("D41").Value = '1.294'
$ws.Range("E41").Value = '  -10.11%  '
$ws.Range("D42").Value = '9.932'
$ws.Range("E42").Value = '  -9.17%  '
$ws.Range("D43").Value = '1.035'
$ws.Range("E43").Value = '  -11.66%  '
$ws.Range("D44").Value = '3.458'
$ws.Range("E44").Value = '  -6.64%  '
$ws.Range("D45").Value = '0.5148'
$ws.Range("E45").Value = '  -10.44%  '
$ws.Range("D46").Value = '11.65'
$ws.Range("E46").Value = '  -10.47%  '
$ws.Range("D47").Value = '0.4908'
$ws.Range("E47").Value = '  -10.40%  '
$ws.Range("D48").Value = '107.44'
$ws.Range("E48").Value = '  -7.68%  '
$ws.Range("D49").Value = '1.706'
$ws.Range("E49").Value = '  -8.51%  '
$ws.Range("D50").Value = '1.007'
$ws.Range("E50").Value = '  +0.65%  '
$ws.Range("D51").Value = '1.024'
$ws.Range("E51").Value = '  -12.06%  '

# Restore default (unstyled) cell style now that values are committed as text
$ws.Range("B2:E51").Style = "Normal"
